$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 86, shifting existing rows 86..173 down to 87..174
$ws.Rows("86:86").Insert()

# Populate the newly inserted row 86 with the new record's data.
$ws.Cells.Item(86, 1).Value = 10
$ws.Cells.Item(86, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(86, 3).Value = "La Araucanía"
$ws.Cells.Item(86, 4).Value = 44789
$ws.Cells.Item(86, 5).Value = 9
$ws.Cells.Item(86, 6).Value = 100112012
$ws.Cells.Item(86, 7).Value = "Espinaca"
$ws.Cells.Item(86, 8).Value = "Sin especificar"
$ws.Cells.Item(86, 9).Value = "Primera"
$ws.Cells.Item(86, 10).Value = 85
$ws.Cells.Item(86, 11).Value = 13000
$ws.Cells.Item(86, 12).Value = 13000
$ws.Cells.Item(86, 13).Value = 13000
$ws.Cells.Item(86, 14).Value = "$/docena de atados"
$ws.Cells.Item(86, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(86, 16).Value = 4333
$ws.Cells.Item(86, 17).Value = 3
$ws.Cells.Item(86, 18).Value = "Hortaliza"
